# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 86
$ws1.Range("F3").Value = 341
$ws1.Range("F4").Value = 4660
$ws1.Range("F5").Value = 45
$ws1.Range("F6").Value = 471

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 86
$ws4.Range("F3").Value = 341
$ws4.Range("F4").Value = 4660
$ws4.Range("F7").Value = 45
$ws4.Range("F8").Value = 471
